$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B6").Value = "internet"
$ws.Range("C6").Value = 123

$ws.Range("C2").Value = 452.25
$ws.Range("C3").Value = 92

$ws.Range("F12").Select()
